$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45874
$ws.Range("B2").Value = 140.7
$ws.Range("C2").Value = 107.94
$ws.Range("D2").Value = 95.40000000000001
$ws.Range("E2").Value = 95.40000000000001
$ws.Range("F2").Value = 94.5
$ws.Range("G2").Value = 91.25
$ws.Range("H2").Value = 101.42
$ws.Range("I2").Value = 106.22
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 30.1
$ws.Range("L2").Value = 18.1
$ws.Range("M2").Value = 1.72
$ws.Range("N2").Value = 0.8
$ws.Range("O2").Value = 0.85
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1.01
$ws.Range("R2").Value = 1.72
$ws.Range("S2").Value = 0.66
$ws.Range("T2").Value = 18.1
$ws.Range("U2").Value = 56
$ws.Range("V2").Value = 84.20999999999999
$ws.Range("W2").Value = 102.62
$ws.Range("X2").Value = 98.36
$ws.Range("Y2").Value = 83.16
$ws.Range("Z2").Value = 58.8
$ws.Range("AA2").Value = "0h-4h"
$ws.Range("AB2").Value = 109.86
$ws.Range("AC2").Value = "0h-2h"
$ws.Range("AD2").Value = 124.32
$ws.Range("AE2").Value = "6h-8h"
$ws.Range("AF2").Value = 103.82
$ws.Range("AG2").Value = "9h-19h"
